# RunControl: make r.yos effective and more analysis on demographics
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RunControl")

# --- 1. Flip the "run it" checkbox (column C) off for the LA-CERA-43 block (rows 14-24) ---
for ($r = 14; $r -le 24; $r++) {
  $ws.Cells.Item($r, 3).Value2 = $False
}

# --- 2. Add two new demographic scenario rows (26 & 27), cloned from the
#        D1F075-mature1 / D1F075-mature1_gn1 rows (19 & 20) but not yet
#        flagged to run (column D = FALSE), for further analysis ---
function Set-ScenarioRow($r, $name, $oVal) {
  $ws.Cells.Item($r, 1).Value2 = $name
  $ws.Cells.Item($r, 2).Value2 = "75% initial Funding; Full smoothing"
  $ws.Cells.Item($r, 3).Value2 = $True
  $ws.Cells.Item($r, 4).Value2 = $False
  $ws.Cells.Item($r, 5).Value2 = "LA-CERA-43.fillin.yos"
  $ws.Cells.Item($r, 6).Value2 = "LA-CERA-43.fillin"
  $ws.Cells.Item($r, 7).Value2 = 1000
  $ws.Cells.Item($r, 8).Value2 = 600
  $ws.Cells.Item($r, 9).Value2 = "average"
  $ws.Cells.Item($r, 10).Value2 = "average"
  $ws.Cells.Item($r, 11).Value2 = "LA-CERA-43.yos"
  $ws.Cells.Item($r, 12).Value2 = "rp2014.hybrid"
  $ws.Cells.Item($r, 13).Value2 = "term.average"
  $ws.Cells.Item($r, 14).Value2 = "nr60er50"
  $ws.Cells.Item($r, 15).Value2 = $oVal
  $ws.Cells.Item($r, 16).Value2 = $False
  $ws.Cells.Item($r, 17).Value2 = 0.022
  $ws.Cells.Item($r, 18).Value2 = 3
  $ws.Cells.Item($r, 19).Value2 = 75
  $ws.Cells.Item($r, 20).Value2 = 50
  $ws.Cells.Item($r, 21).Value2 = 60
  $ws.Cells.Item($r, 22).Value2 = 0.02
  $ws.Cells.Item($r, 23).Value2 = 0
  $ws.Cells.Item($r, 24).Value2 = 10
  $ws.Cells.Item($r, 25).Value2 = 0.04
  $ws.Cells.Item($r, 26).Value2 = 0.04
  $ws.Cells.Item($r, 27).Value2 = 0.03
  $ws.Cells.Item($r, 28).Value2 = 0.01
  $ws.Cells.Item($r, 29).Value2 = 0.075
  $ws.Cells.Item($r, 30).Value2 = "simple"
  $ws.Cells.Item($r, 31).Value2 = 0.0822
  $ws.Cells.Item($r, 32).Value2 = 0.12
  $ws.Cells.Item($r, 33).Value2 = "EAN.CP"
  $ws.Cells.Item($r, 34).Value2 = "open"
  $ws.Cells.Item($r, 35).Value2 = "cp"
  $ws.Cells.Item($r, 36).Value2 = 30
  $ws.Cells.Item($r, 37).Value2 = "method1"
  $ws.Cells.Item($r, 38).Value2 = 5
  $ws.Cells.Item($r, 39).Value2 = 200
  $ws.Cells.Item($r, 40).Value2 = "MA"
  $ws.Cells.Item($r, 41).Value2 = 1
  $ws.Cells.Item($r, 42).Value2 = "AL_pct"
  $ws.Cells.Item($r, 43).Value2 = 0.75
  $ws.Cells.Item($r, 44).Value2 = 200
  $ws.Cells.Item($r, 45).Value2 = "ADC"
  $ws.Cells.Item($r, 46).Value2 = 0.25
  $ws.Cells.Item($r, 47).Value2 = 0.145
  $ws.Cells.Item($r, 48).Value2 = 0.05
  $ws.Cells.Item($r, 49).Value2 = $False
  $ws.Cells.Item($r, 50).Value2 = $True
  $ws.Cells.Item($r, 51).Value2 = $False

  # Scenario-name cell (column A) carries the green "new scenario" highlight
  $ws.Cells.Item($r, 1).Interior.Color = 5296274
}

Set-ScenarioRow 26 "D1F075-mature1_lowB" 0
Set-ScenarioRow 27 "D1F075-mature1_gn1_lowB" -0.01

# --- 3. Scroll / selection bookkeeping so the view lands on the new rows ---
$ws.Range("D29:D30").Select()
